$wb = $excel.ActiveWorkbook

$wsBabies = $wb.Worksheets.Item("babies")
$wsRooms  = $wb.Worksheets.Item("rooms")

# Add the new "treatment" column to the babies sheet (D1)
$wsBabies.Range("D1").Value = "treatment"

# Add the new "treatment" column to the rooms sheet (I1)
$wsRooms.Range("I1").Value = "treatment"

# Update selections / active cells to match the edited columns
$wsBabies.Range("D1").Select()
$wsRooms.Range("I1").Select()

# Make "babies" sheet the active tab
$wsBabies.Activate()
